$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) values to latest scrape.
# NumberFormat is temporarily set to Text ("@") before assigning numeric-looking
# strings so Excel keeps them as text (matching original inlineStr cells), then
# the cell style is reset back to Normal so no stray formatting is introduced.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "67.656.17"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +1.40%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.613.06"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("E4").Value = "  +0.05%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "602.33"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.62%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "154.30"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.57%  "

$ws.Range("E7").Value = "  +0.04%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.612.00"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.72%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.126"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +8.90%  "

$ws.Range("E11").Value = "  +0.64%  "

$ws.Range("E12").Value = "  +0.93%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.354"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.60%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "28.01"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.00%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.0000186"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +3.24%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.093.45"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.93%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "67.527.84"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.30%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.614.79"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.81%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "11.26"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -0.50%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "365.55"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +3.50%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "7.59"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.59%  "

$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("E23").Value = "  +3.54%  "

$ws.Range("E24").Value = "  +0.00%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "70.01"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.12%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "10.14"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -2.74%  "

$ws.Range("E27").Value = "  +2.15%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.746.21"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.66%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "582.19"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.82%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.31%  "

$ws.Range("E31").Value = "  -0.99%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "7.93"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.32%  "

$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("E34").Value = "  -2.34%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("E36").Value = "  -1.85%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "4.97"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.65%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "19.41"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "155.35"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +1.01%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.372"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.77%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "5.40"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.75%  "

$ws.Range("E42").Value = "  +2.38%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.64"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +2.35%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "41.11"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.77%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "16.43"
$cell.Style = "Normal"

$ws.Range("E46").Value = "  +0.03%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "156.57"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.06%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.0₆0287"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -6.82%  "

$ws.Range("E49").Value = "  -0.24%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "21.02"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.75%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.623"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.80%  "
